$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("12/9/2022", "-3", "store 1", "Ximena Leyva", "-"),
    @("9/27/2022", "-43.78", "another store", "Ximena Leyva", "-"),
    @("10/27/2022", "-142", "shoe store", "Ximena Leyva", "shoes"),
    @("11/8/2022", "-34", "make up store", "Ximena Leyva", "-"),
    @("12/3/2022", "-45", "store1", "Ximena Leyva", "-"),
    @("10/5/2022", "+430", "funding 2", "Ximena Leyva", "-"),
    @("11/10/2022", "+534", "source 4", "Ximena Leyva", "-"),
    @("9/9/2022", "+24", "donation", "Ximena Leyva", "-"),
    @("11/1/2022", "+1000", "CPA", "Ximena Leyva", "-"),
    @("12/18/2022", "+2500", "massive donation", "Ximena Leyva", "-"),
    @("12/11/2022", "+430", "cpa", "Ximena Leyva", "-")
)

$ws.Range("A2:E12").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
